$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold numeric-looking text (Protocolo / CPF) must keep their
# original "text" cell type, so force a text number format before writing
# the value, then revert the cell style so no stray formatting is left
# behind.
$textCells = @("E2", "E3", "G2", "G3")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Operadora (same new value on both rows)
$ws.Range("A2").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("A3").Value = "417823 - PREMIUM SAÚDE S.A"

# Row 2 (Notificação, Demanda, Protocolo, Beneficiário, CPF, Descrição, Prazo)
$ws.Range("C2").Value = 45027.36063657407
$ws.Range("D2").Value = 12163407
$ws.Range("E2").Value = "8600473"
$ws.Range("F2").Value = "AYLA ALVES COELHO"
$ws.Range("G2").Value = "19014458606"
$ws.Range("H2").Value = "minha filha usa NEOCATE desde os 4 meses de idade, foi solicitado o pedido de intolerancia a lactose no dia 22/3/2023 porem nao tenho retorno quando ligo dizem, que pode levar ate 21 dias uteis mas em consulta ao site da ANS esse prazo seria pra procedimentos PAC, que nao e o caso de um exame de sangue para detectar intolerância a lactose, preciso de um retorno visto que o prazo de 10 dias uteis finalizou em 06/04/2023."
$ws.Range("I2").Value = 6

# Row 3 (Notificação, Demanda, Protocolo, Beneficiário, CPF, Descrição, Prazo)
$ws.Range("C3").Value = 45027.46474537037
$ws.Range("D3").Value = 12163869
$ws.Range("E3").Value = "8601052"
$ws.Range("F3").Value = "PRISCILA APARECIDA SANTOS FRANCISCO"
$ws.Range("G3").Value = "10136083641"
$ws.Range("H3").Value = "Solicitei com 10 dias de antecedência a marcação do exame de ultrassonografia endovaginal,na clínica Santa Helena ltda através do plano,ao chegar no estabelecimento no dia do exame a operadora do plano negou meu procedimento alegando que a clinica não estava mais cadastrada,porém não recebi nenhum contanto prévio mesmo estando com o agendamento feito a dias."
$ws.Range("I3").Value = 6

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
